$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 499.5
$ws.Range("I20").Value = 499.5
$ws.Range("K20").Value = 499.5
$ws.Range("M20").Value = -269.5
$ws.Range("H35").Value = 499.5
$ws.Range("I35").Value = 499.5
$ws.Range("K35").Value = 499.5
$ws.Range("M35").Value = -120.5
$ws.Range("H112").Value = 0
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("M112").ClearContents()
$ws.Range("N112").ClearContents()
$ws.Range("H137").Value = 2190.3096
$ws.Range("I137").Value = 1591.0435
$ws.Range("J137").Value = 2915.7368
$ws.Range("K137").Value = 4773.1305
$ws.Range("L137").Value = 8747.2104
$ws.Range("M137").Value = -2223.1305
$ws.Range("N137").Value = -13847.2104
$ws.Range("H138").Value = 10093.765
$ws.Range("J138").Value = 10866.417
$ws.Range("L138").Value = 32599.251
$ws.Range("N138").Value = -42879.251

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 467.8
$ws.Range("I5").Value = 434.75
$ws.Range("J5").Value = 600
$ws.Range("K5").Value = 434.75
$ws.Range("L5").Value = 600
$ws.Range("M5").Value = -322.75
$ws.Range("N5").Value = -824
$ws.Range("H32").Value = 1104.037
$ws.Range("I32").Value = 629.95917
$ws.Range("K32").Value = 629.95917
$ws.Range("M32").Value = -342.95917
$ws.Range("H61").Value = 6349.9033
$ws.Range("I61").Value = 4709
$ws.Range("J61").Value = 10361
$ws.Range("K61").Value = 4709
$ws.Range("L61").Value = 10361
$ws.Range("M61").Value = -4497
$ws.Range("N61").Value = -10785
$ws.Range("H74").Value = 2661.5356
$ws.Range("I74").Value = 1785.238
$ws.Range("J74").Value = 5290.4287
$ws.Range("K74").Value = 1785.238
$ws.Range("L74").Value = 5290.4287
$ws.Range("M74").Value = -911.2380000000001
$ws.Range("N74").Value = -7038.4287
$ws.Range("H77").Value = 2661.5356
$ws.Range("I77").Value = 1785.238
$ws.Range("J77").Value = 5290.4287
$ws.Range("K77").Value = 8926.190000000001
$ws.Range("L77").Value = 26452.1435
$ws.Range("M77").Value = -4558.190000000001
$ws.Range("N77").Value = -35188.14350000001
$ws.Range("H104").Value = 39999
$ws.Range("J104").Value = 39999
$ws.Range("L104").Value = 39999
$ws.Range("N104").Value = -46987
$ws.Range("H132").Value = 3861.1143
$ws.Range("I132").Value = 3562.6128
$ws.Range("J132").Value = 6174.5
$ws.Range("K132").Value = 10687.8384
$ws.Range("L132").Value = 18523.5
$ws.Range("M132").Value = -8157.838400000001
$ws.Range("N132").Value = -23583.5
$ws.Range("H136").Value = 6349.9033
$ws.Range("I136").Value = 4709
$ws.Range("J136").Value = 10361
$ws.Range("K136").Value = 14127
$ws.Range("L136").Value = 31083
$ws.Range("M136").Value = -11577
$ws.Range("N136").Value = -36183

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 467.8
$ws.Range("I4").Value = 434.75
$ws.Range("J4").Value = 600
$ws.Range("K4").Value = 434.75
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = -319.75
$ws.Range("N4").Value = -830
$ws.Range("H107").Value = 1375.5555
$ws.Range("I107").Value = 1547.2142
$ws.Range("J107").Value = 774.75
$ws.Range("K107").Value = 1547.2142
$ws.Range("L107").Value = 774.75
$ws.Range("M107").Value = 372.7858000000001
$ws.Range("N107").Value = -4614.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 94.666664
$ws.Range("I7").Value = 94.666664
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 94.666664
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 18.333336
$ws.Range("N7").ClearContents()
$ws.Range("H22").Value = 875.5454999999999
$ws.Range("I22").Value = 483.83334
$ws.Range("K22").Value = 483.83334
$ws.Range("M22").Value = -133.83334
$ws.Range("H31").Value = 3788.6924
$ws.Range("I31").Value = 3161.4546
$ws.Range("J31").Value = 3956.9756
$ws.Range("K31").Value = 3161.4546
$ws.Range("L31").Value = 3956.9756
$ws.Range("M31").Value = -2866.4546
$ws.Range("N31").Value = -4546.9756
$ws.Range("H34").Value = 3788.6924
$ws.Range("I34").Value = 3161.4546
$ws.Range("J34").Value = 3956.9756
$ws.Range("K34").Value = 3161.4546
$ws.Range("L34").Value = 3956.9756
$ws.Range("M34").Value = -2959.4546
$ws.Range("N34").Value = -4360.9756

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1014
$ws.Range("I2").Value = 1346
$ws.Range("J2").Value = 18
$ws.Range("K2").Value = 8076
$ws.Range("L2").Value = 108
$ws.Range("M2").Value = -7963
$ws.Range("N2").Value = -334
$ws.Range("H5").Value = 2664.0454
$ws.Range("I5").Value = 1866.7778
$ws.Range("J5").Value = 3216
$ws.Range("K5").Value = 5600.3334
$ws.Range("L5").Value = 9648
$ws.Range("M5").Value = -5488.3334
$ws.Range("N5").Value = -9872
$ws.Range("H58").Value = 8335.333000000001
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 8335.333000000001
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 25005.999
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -25261.999
$ws.Range("H68").Value = 3422.5293
$ws.Range("J68").Value = 3645.5334
$ws.Range("L68").Value = 10936.6002
$ws.Range("N68").Value = -12558.6002
$ws.Range("H71").Value = 3422.5293
$ws.Range("J71").Value = 3645.5334
$ws.Range("L71").Value = 32809.8006
$ws.Range("N71").Value = -40921.8006
$ws.Range("H107").Value = 1986.8695
$ws.Range("I107").Value = 259
$ws.Range("J107").Value = 2065.4092
$ws.Range("K107").Value = 777
$ws.Range("L107").Value = 6196.2276
$ws.Range("M107").Value = 1143
$ws.Range("N107").Value = -10036.2276
$ws.Range("H135").Value = 2664.0454
$ws.Range("I135").Value = 1866.7778
$ws.Range("J135").Value = 3216
$ws.Range("K135").Value = 16801.0002
$ws.Range("L135").Value = 28944
$ws.Range("M135").Value = -14266.0002
$ws.Range("N135").Value = -34014

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3623.75
$ws.Range("I102").Value = 2297.2856
$ws.Range("J102").Value = 4655.4443
$ws.Range("K102").Value = 2297.2856
$ws.Range("L102").Value = 4655.4443
$ws.Range("M102").Value = -675.2856000000002
$ws.Range("N102").Value = -7899.4443

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2439
$ws.Range("I40").Value = 2170
$ws.Range("J40").Value = 3138.4
$ws.Range("K40").Value = 2170
$ws.Range("L40").Value = 3138.4
$ws.Range("M40").Value = -2034
$ws.Range("N40").Value = -3410.4
$ws.Range("H63").Value = 67488
$ws.Range("I63").Value = 69977
$ws.Range("J63").Value = 64999
$ws.Range("K63").Value = 69977
$ws.Range("L63").Value = 64999
$ws.Range("M63").Value = -69228
$ws.Range("N63").Value = -66497
$ws.Range("H66").Value = 67488
$ws.Range("I66").Value = 69977
$ws.Range("J66").Value = 64999
$ws.Range("K66").Value = 209931
$ws.Range("L66").Value = 194997
$ws.Range("M66").Value = -206187
$ws.Range("N66").Value = -202485
$ws.Range("H122").Value = 2180.087
$ws.Range("I122").Value = 2010.1538
$ws.Range("K122").Value = 6030.4614
$ws.Range("M122").Value = -3580.4614

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H104").Value = 16833.2
$ws.Range("J104").Value = 16833.2
$ws.Range("L104").Value = 16833.2
$ws.Range("N104").Value = -23821.2
$ws.Range("H140").Value = 89489.71000000001
$ws.Range("J140").Value = 89489.71000000001
$ws.Range("L140").Value = 89489.71000000001
$ws.Range("N140").Value = -99849.71000000001
